$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header values in row 1: P1=14, Q1=15 (same bold/border style as O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# For rows 2 through 25, swap values in columns I/K and M/O, and add new columns P and Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2
}
